$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextCell 2 4 '27.868.33'
$ws.Cells.Item(2, 5).Value = '  -0.37%  '
Set-TextCell 3 4 '1.624.05'
$ws.Cells.Item(3, 5).Value = '  -1.03%  '
$ws.Cells.Item(4, 5).Value = '  -0.18%  '
Set-TextCell 5 4 '211.09'
$ws.Cells.Item(5, 5).Value = '  -1.10%  '
$ws.Cells.Item(6, 5).Value = '  -0.34%  '
$ws.Cells.Item(7, 5).Value = '  -0.19%  '
Set-TextCell 8 4 '23.44'
$ws.Cells.Item(8, 5).Value = '  -0.82%  '
$ws.Cells.Item(9, 5).Value = '  -2.08%  '
Set-TextCell 10 4 '0.0611'
$ws.Cells.Item(10, 5).Value = '  -0.59%  '
$ws.Cells.Item(11, 5).Value = '  +0.00%  '
Set-TextCell 12 4 '1.855.76'
Set-TextCell 13 4 '1.617.56'
$ws.Cells.Item(13, 5).Value = '  -1.44%  '
$ws.Cells.Item(14, 5).Value = '  -1.96%  '
$ws.Cells.Item(15, 5).Value = '  -2.36%  '
Set-TextCell 16 4 '65.35'
$ws.Cells.Item(16, 5).Value = '  -0.86%  '
Set-TextCell 17 4 '27.852.60'
$ws.Cells.Item(17, 5).Value = '  -0.47%  '
Set-TextCell 18 4 '229.66'
$ws.Cells.Item(18, 5).Value = '  -1.29%  '
$ws.Cells.Item(19, 5).Value = '  +0.67%  '
Set-TextCell 20 4 '0.0₃0721'
$ws.Cells.Item(20, 5).Value = '  -0.41%  '
$ws.Cells.Item(21, 5).Value = '  -0.22%  '
$ws.Cells.Item(22, 5).Value = '  -0.94%  '
Set-TextCell 23 4 '10.14'
$ws.Cells.Item(23, 5).Value = '  -5.84%  '
$ws.Cells.Item(24, 5).Value = '  -2.63%  '
Set-TextCell 25 4 '154.66'
$ws.Cells.Item(25, 5).Value = '  +1.95%  '
$ws.Cells.Item(26, 5).Value = '  -1.38%  '
$ws.Cells.Item(27, 5).Value = '  -0.11%  '
$ws.Cells.Item(28, 5).Value = '  -1.37%  '
$ws.Cells.Item(29, 5).Value = '  -0.15%  '
$ws.Cells.Item(30, 5).Value = '  -0.79%  '
$ws.Cells.Item(32, 5).Value = '  +2.02%  '
$ws.Cells.Item(33, 2).Value = 'Maker'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextCell 33 4 '1.400.29'
$ws.Cells.Item(33, 5).Value = '  -0.30%  '
$ws.Cells.Item(34, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextCell 34 4 '3.07'
$ws.Cells.Item(34, 5).Value = '  -1.12%  '
$ws.Cells.Item(35, 5).Value = '  -0.42%  '
Set-TextCell 36 4 '0.998'
$ws.Cells.Item(36, 5).Value = '  +9.02%  '
$ws.Cells.Item(37, 5).Value = '  -1.35%  '
$ws.Cells.Item(38, 5).Value = '  +0.99%  '
Set-TextCell 39 4 '0.554'
$ws.Cells.Item(39, 5).Value = '  -0.54%  '
Set-TextCell 40 4 '0.859'
$ws.Cells.Item(40, 5).Value = '  -3.06%  '
$ws.Cells.Item(41, 5).Value = '  -0.09%  '
$ws.Cells.Item(42, 5).Value = '  -0.26%  '
$ws.Cells.Item(43, 2).Value = 'Aave'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextCell 43 4 '65.89'
$ws.Cells.Item(43, 5).Value = '  -0.68%  '
Set-TextCell 44 4 '5.48'
$ws.Cells.Item(44, 5).Value = '  -0.17%  '
$ws.Cells.Item(45, 2).Value = 'RenderToken'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell 45 4 '1.83'
$ws.Cells.Item(45, 5).Value = '  -2.87%  '
$ws.Cells.Item(46, 5).Value = '  -0.83%  '
Set-TextCell 47 4 '1.765.76'
$ws.Cells.Item(47, 5).Value = '  -0.90%  '
Set-TextCell 48 4 '87.94'
$ws.Cells.Item(48, 5).Value = '  -0.12%  '
Set-TextCell 49 4 '0.102'
$ws.Cells.Item(49, 5).Value = '  +1.83%  '
$ws.Cells.Item(50, 2).Value = 'BabyDogeCoin'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextCell 50 4 '0.0₆0103'
$ws.Cells.Item(50, 5).Value = '  -2.82%  '
$ws.Cells.Item(51, 2).Value = 'Cronos'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell 51 4 '0.0503'
$ws.Cells.Item(51, 5).Value = '  -0.54%  '
